# Apply the "IBM -> LinuxForHealth" rebrand + version bump edit described
# by the commit:
#   Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500...

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Metadata": top-level StructureDefinition property/value table
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-item-detail-classification"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------------
# Sheet "Elements": element definition table
# ---------------------------------------------------------------------
$elem = $wb.Worksheets.Item("Elements")

# Row 2 = "Extension" root element: its Constraint(s) cell incorrectly
# carried the ele-1/ext-1 constraint text (that belongs on the
# Extension.extension row) - clear it.
$elem.Range("AI2").Value = ""

# Row 5 = "Extension.url": Fixed Value is the StructureDefinition's own
# canonical URL.
$elem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-item-detail-classification"

# Row 7 = "Extension.value[x]" / valueCodeableConcept slice: Binding
# Value Set points at the companion ValueSet, also rebranded.
$elem.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/claim-item-detail-classification"

# Column Y ("Binding Value Set") widens to fit the longer URL text.
$elem.Columns.Item(25).ColumnWidth = 68
